# Add Q4-2022 data: a new "2022-Q4" worksheet (copied from the "2022-Q3"
# sheet so it inherits the same layout/styling) placed right after the
# "总计" (totals) sheet, and a new leading row on the "总计" sheet summarizing
# the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" (same columns/
#    header/styles) and inserting it immediately before "2022-Q3".
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3, $null)
$q4 = $wb.ActiveSheet
$q4.Name = "2022-Q4"

# Wipe the copied (2022-Q3) data so we can write 2022-Q4's own rows.
$q4.Range("B2:H10").ClearContents()
$q4.Range("B2:H10").ClearFormats()
$q4.Range("A2:A10").ClearContents()

# Columns B-G hold text (fund code / name / percentages as strings, not
# numbers) in every quarter sheet - force text storage before assigning so
# e.g. "008115" doesn't get coerced into a number and lose its leading zero.
$q4.Range("B2:G6").NumberFormat = "@"

$q4rows = @(
    @("008115", "天弘中证红利低波动100指数C",               "2.67", "94.95", "1.75", "0.0467", 7),
    @("515100", "景顺长城中证红利低波动100ETF",               "1.96", "99.17", "1.83", "0.0359", 7),
    @("008114", "天弘中证红利低波动100指数A",               "1.98", "94.95", "1.75", "0.0346", 7),
    @("562530", "华夏中证智选1000价值稳健策略ETF",             "0.36", "96.22", "0.92", "0.0033", 8),
    @("005770", "信澳中证沪港深高股息精选指数",               "0.13", "23.47", "0.73", "0.0009", 1)
)

for ($i = 0; $i -lt $q4rows.Length; $i++) {
    $r = $i + 2
    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $q4rows[$i][0]
    $q4.Cells.Item($r, 3).Value = $q4rows[$i][1]
    $q4.Cells.Item($r, 4).Value = $q4rows[$i][2]
    $q4.Cells.Item($r, 5).Value = $q4rows[$i][3]
    $q4.Cells.Item($r, 6).Value = $q4rows[$i][4]
    $q4.Cells.Item($r, 7).Value = $q4rows[$i][5]
    $q4.Cells.Item($r, 8).Value = $q4rows[$i][6]
}

# Column A (row index, 0-based) keeps the same bordered/centered style as the
# header row uses for the rest of column A - re-stamp it across the rows we
# just (re)populated so every data row matches (the copy only carried it for
# the two rows that previously had data).
$q4.Range("A2").Copy()
$q4.Range("A3:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new leading row for 2022-Q4
#    and keep the running index in column A sequential.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows("2:2").Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.12

for ($r = 3; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# Insert() copies the header row's (bold/bordered) formatting into the new
# row - restore the plain data-row look by copying formats down from the row
# right below (untouched by the insert).
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Leave the selection/active tab on "2020-Q4" (the last sheet), matching
#    the original workbook's state.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item("2020-Q4")
$lastSheet.Activate()
$null = $lastSheet.Range("A1").Select()
